# Applies the "adding averages and more checks" update:
#  - Training Dashboard: PERIOD TO EXPIRE (H) shrinks by 8 and LAST UPDATE (I)
#    moves from 08-Sep-2025 to 16-Sep-2025 for every data row (3-29)
#  - Exam Dashboard: COMMENTS column (E) text normalised to "date is valid",
#    its column narrowed, and the header/title styling is tweaked so the
#    bold header text is white instead of black
#  - the oversized "title" font (bold, 14pt) is dropped in favour of the
#    regular bold font so both the sheet titles and the table headers share
#    the same (now white) bold font

$wb = $excel.ActiveWorkbook

$trainingWs = $wb.Worksheets.Item("Training Dashboard")
$examWs     = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------
# 1) Training Dashboard - update PERIOD TO EXPIRE / LAST UPDATE columns
# ---------------------------------------------------------------------
$periodToExpire = @{
    3  = 241
    4  = 213
    5  = 239
    6  = 363
    7  = 244
    8  = 336
    9  = 219
    10 = 247
    11 = 238
    12 = 245
    13 = 223
    14 = 349
    15 = 665
    16 = 258
    17 = 413
    18 = 413
    19 = -23
    20 = -83
    21 = -106
    22 = -34
    23 = -34
    24 = 268
    25 = 313
    26 = 313
    27 = 313
    28 = 323
    29 = 348
}

# Make sure the "LAST UPDATE" cells keep their new date as literal text
# (matching the original inline-string content) instead of being
# auto-converted into a date serial number.
$trainingWs.Range("I3:I29").NumberFormat = "@"

foreach ($row in 3..29) {
    $trainingWs.Cells.Item($row, 8).Value = $periodToExpire[$row]
    $trainingWs.Cells.Item($row, 9).Value = "16-Sep-2025"
}

# ---------------------------------------------------------------------
# 2) Exam Dashboard - normalise COMMENTS column and narrow it
# ---------------------------------------------------------------------
foreach ($row in 3..12) {
    $examWs.Cells.Item($row, 5).Value = "date is valid"
}

# Column E was 44 characters wide, now 15. Excel's ColumnWidth property is
# offset from the raw column width stored in the file by 5/6 of a
# character, so compensate to land exactly on 15.
$examWs.Columns("E").ColumnWidth = 15 - 5/6

# ---------------------------------------------------------------------
# 3) Styling - drop the big 14pt title font; both the dashboard titles
#    and the table headers now use the same bold, white font.
# ---------------------------------------------------------------------
$trainingWs.Range("A1:K1").Font.Bold = $true
$trainingWs.Range("A1:K1").Font.Size = 11
$trainingWs.Range("A1:K1").Font.Color = 16777215

$trainingWs.Range("A2:K2").Font.Bold = $true
$trainingWs.Range("A2:K2").Font.Color = 16777215

$examWs.Range("A1:G1").Font.Bold = $true
$examWs.Range("A1:G1").Font.Size = 11
$examWs.Range("A1:G1").Font.Color = 16777215

$examWs.Range("A2:G2").Font.Bold = $true
$examWs.Range("A2:G2").Font.Color = 16777215
